$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that were removed from the dataset (by original row
# number, highest first so row numbers of the other target row don't shift
# before we get to it).
$ws.Rows.Item(28).Delete()   # "SC 92"
$ws.Rows.Item(26).Delete()   # "RM 232"

# After the deletions, rows 27.. and onward shifted up by two, and row 28
# ("SC 92") also already removed, so the surviving rows now sit exactly
# where the target sheet expects them (header + 32 data rows => A1:F33).

# --- Scattered value fixes on the still-present rows (not touched by the
# row deletions above) ---

# RM 2 (row 2): E2 value removed -> now blank
$ws.Range("E2").Value = $null

# RM 14 (row 5): E5 blank -> now has a value
$ws.Range("E5").Value = -5

# RM 21 (row 6): D6 and E6 blank -> now have values
$ws.Range("D6").Value = -14.2
$ws.Range("E6").Value = -5.7

# RM 38 (row 8): D8 value removed -> now blank
$ws.Range("D8").Value = $null

# RM 42 (row 9): E9 value removed -> now blank
$ws.Range("E9").Value = $null

# RM 52 a (row 10): E10 value removed -> now blank
$ws.Range("E10").Value = $null

# RM 81 (row 12): D12 blank -> now has a value
$ws.Range("D12").Value = -14.1

# RM 90 (row 14): D14 value removed -> now blank
$ws.Range("D14").Value = $null

# RM 116 (row 17): D17 blank -> now has a value
$ws.Range("D17").Value = -14.7

# RM 120 (row 18): D18 blank -> now has a value
$ws.Range("D18").Value = -15.2

# RM 125 (row 19): D19 value removed -> now blank
$ws.Range("D19").Value = $null

# RM 134 (row 20): D20 value removed -> now blank
$ws.Range("D20").Value = $null

# RM 140 (row 23): D23 blank -> now has a value
$ws.Range("D23").Value = -13.9

# RM 142a (row 24): E24 blank -> now has a value
$ws.Range("E24").Value = -8.1

# --- Rows that used to be 27-35 (SC 5 .. SC 232), now 26-33 after the two
# row deletions above. A handful of their cells differ from the simple
# "shift up" copy, so correct those specific cells to match the target. ---

# SC 101 (now row 27): B27 blank -> now has a value; D27 value removed -> now blank
$ws.Range("B27").Value = -20.4
$ws.Range("D27").Value = $null

# SC 105 (now row 28): B28 and E28 values removed -> now blank
$ws.Range("B28").Value = $null
$ws.Range("E28").Value = $null

# SC 119 (now row 29): B29 value removed -> now blank
$ws.Range("B29").Value = $null

# SC 120 (now row 30): B30 and E30 blank -> now have values
$ws.Range("B30").Value = -19.7
$ws.Range("E30").Value = -5.7

# SC 193 (now row 32): B32 value removed -> now blank
$ws.Range("B32").Value = $null
